$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (K2:T2) with new TPM-derived values
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.9242423333333334
$ws.Range("N2").Value = 2.772727
$ws.Range("O2").Value = 0.7379131210038523
$ws.Range("P2").Value = 0.8085498937233963
$ws.Range("Q2").Value = 0.1401619660115556
$ws.Range("R2").Value = 1.261457694104
$ws.Range("S2").Value = 0.7379131210038523
$ws.Range("T2").Value = 0.8085498937233963

# Add new row 3 for the MuSCs target cluster
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Kng1"
$ws.Range("C3").Value = "Bdkrb2"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1516506666666667
$ws.Range("H3").Value = 0.454952
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.328266
$ws.Range("N3").Value = 0.656532
$ws.Range("O3").Value = 0.2620868789961477
$ws.Range("P3").Value = 0.1914501062766038
$ws.Range("Q3").Value = 0.049781757744
$ws.Range("R3").Value = 0.298690546464
$ws.Range("S3").Value = 0.2620868789961477
$ws.Range("T3").Value = 0.1914501062766038
